$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "62.313.93"
$ws.Cells.Item(2, 5).Value = "  -2.30%  "

$ws.Cells.Item(3, 4).Value = "2.638.76"
$ws.Cells.Item(3, 5).Value = "  -3.42%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "548.76"

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "153.94"
$ws.Cells.Item(6, 5).Value = "  -4.69%  "

$ws.Cells.Item(7, 5).Value = "  +0.08%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.585"
$ws.Cells.Item(8, 5).Value = "  -1.89%  "

$ws.Cells.Item(9, 5).Value = "  -4.61%  "

$ws.Cells.Item(10, 5).Value = "  -4.33%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "5.41"
$ws.Cells.Item(11, 5).Value = "  -3.73%  "

$ws.Cells.Item(12, 5).Value = "  -4.89%  "

$ws.Cells.Item(13, 4).Value = "3.106.33"
$ws.Cells.Item(13, 5).Value = "  -3.52%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "25.59"
$ws.Cells.Item(14, 5).Value = "  -5.02%  "

$ws.Cells.Item(15, 4).Value = "62.237.16"
$ws.Cells.Item(15, 5).Value = "  -2.25%  "

$ws.Cells.Item(16, 5).Value = "  -4.11%  "

$ws.Cells.Item(17, 4).Value = "2.640.74"
$ws.Cells.Item(17, 5).Value = "  -3.80%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "11.66"
$ws.Cells.Item(18, 5).Value = "  -5.61%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "4.51"

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "338.76"
$ws.Cells.Item(20, 5).Value = "  -4.53%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "6.08"
$ws.Cells.Item(21, 5).Value = "  -8.02%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "0.999"
$ws.Cells.Item(22, 5).Value = "  +0.06%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.500"
$ws.Cells.Item(23, 5).Value = "  -3.81%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "62.60"
$ws.Cells.Item(24, 5).Value = "  -2.89%  "

$ws.Cells.Item(25, 5).Value = "  -1.25%  "

$ws.Cells.Item(26, 5).Value = "  +0.07%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "7.98"
$ws.Cells.Item(27, 5).Value = "  -4.71%  "

$ws.Cells.Item(28, 2).Value = "PEPE"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Cells.Item(28, 4).Value = "0.0₃0833"
$ws.Cells.Item(28, 5).Value = "  -8.26%  "

$ws.Cells.Item(29, 2).Value = "Fetch.AI"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "1.36"
$ws.Cells.Item(29, 5).Value = "  +0.55%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "7.05"
$ws.Cells.Item(30, 5).Value = "  -1.45%  "

$ws.Cells.Item(31, 5).Value = "  -5.89%  "

$ws.Cells.Item(32, 2).Value = "USDe"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.999"
$ws.Cells.Item(32, 5).Value = "  +0.05%  "

$ws.Cells.Item(33, 2).Value = "Monero"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "159.81"
$ws.Cells.Item(33, 5).Value = "  -3.93%  "

$ws.Cells.Item(34, 5).Value = "  -3.80%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.42"
$ws.Cells.Item(35, 5).Value = "  -3.67%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "19.15"
$ws.Cells.Item(36, 5).Value = "  -4.69%  "

$ws.Cells.Item(37, 5).Value = "  -4.69%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "332.58"
$ws.Cells.Item(38, 5).Value = "  -3.30%  "

$ws.Cells.Item(39, 5).Value = "  -3.04%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.900"
$ws.Cells.Item(40, 5).Value = "  -7.65%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "37.98"
$ws.Cells.Item(41, 5).Value = "  -1.78%  "

$ws.Cells.Item(42, 5).Value = "  -4.28%  "

$ws.Cells.Item(43, 5).Value = "  +0.06%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "20.33"
$ws.Cells.Item(44, 5).Value = "  -6.65%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.606"
$ws.Cells.Item(45, 5).Value = "  -3.80%  "

$ws.Cells.Item(46, 2).Value = "WhiteBITCoin"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "10.98"
$ws.Cells.Item(46, 5).Value = "  -0.66%  "

$ws.Cells.Item(47, 2).Value = "EnergySwap"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "19.69"
$ws.Cells.Item(47, 5).Value = "  -6.43%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.0545"
$ws.Cells.Item(48, 5).Value = "  -6.75%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.0958"
$ws.Cells.Item(49, 5).Value = "  -4.04%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "127.16"
$ws.Cells.Item(50, 5).Value = "  -3.73%  "

$ws.Cells.Item(51, 5).Value = "  -5.46%  "
